$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Insert a new row above row 4 (shifts existing rows 4-26 down to 5-27,
# and picks up formatting from row 3 as Excel normally does).
$ws.Rows("4:4").Insert()

# Fill in the data for the newly inserted task row.
$ws.Range("B4").Value = "Création de la base de données"
$ws.Range("C4").Value = "Partir les scripts"
$ws.Range("D4").Value = "William Themens"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "En cours"

# Update status of the first two tasks.
$ws.Range("F2").Value = "Terminé"
$ws.Range("F3").Value = "En cours"

# Update hours logged on a couple of existing tasks (now shifted to rows 6 & 8).
$ws.Range("E6").Value = 10
$ws.Range("E8").Value = 4

# Add a total-hours formula.
$ws.Range("E20").Formula = "=SUM(E2:E19)"

# Match the recorded selection after the edits.
$ws.Range("F4").Select() | Out-Null
